# Auto-generated edit script: updates market-price derived columns (H-N)
# in the Leve profit tracking sheets to reflect refreshed Market Board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1340
$ws.Range("I9").Value = 1340
$ws.Range("K9").Value = 1340
$ws.Range("M9").Value = -1171

$ws.Range("H12").Value = 868.6667
$ws.Range("I12").Value = 868.6667
$ws.Range("K12").Value = 868.6667
$ws.Range("M12").Value = -698.6667

$ws.Range("H40").Value = 1329.6
$ws.Range("I40").Value = 1329.6
$ws.Range("K40").Value = 1329.6
$ws.Range("M40").Value = -1154.6

$ws.Range("H53").Value = 53.5
$ws.Range("I53").Value = 58.4
$ws.Range("J53").Value = 29
$ws.Range("K53").Value = 58.4
$ws.Range("L53").Value = 29
$ws.Range("M53").Value = 578.6
$ws.Range("N53").Value = -1303

$ws.Range("H70").Value = 5689.731
$ws.Range("I70").Value = 3646.8333
$ws.Range("J70").Value = 7440.7856
$ws.Range("K70").Value = 10940.4999
$ws.Range("L70").Value = 22322.3568
$ws.Range("M70").Value = -10670.4999
$ws.Range("N70").Value = -22862.3568

$ws.Range("H73").Value = 5689.731
$ws.Range("I73").Value = 3646.8333
$ws.Range("J73").Value = 7440.7856
$ws.Range("K73").Value = 10940.4999
$ws.Range("L73").Value = 22322.3568
$ws.Range("M73").Value = -10004.4999
$ws.Range("N73").Value = -24194.3568

$ws.Range("H80").Value = 628.4091
$ws.Range("I80").Value = 331.27274
$ws.Range("J80").Value = 925.5454999999999
$ws.Range("K80").Value = 993.81822
$ws.Range("L80").Value = 2776.6365
$ws.Range("M80").Value = 4.181780000000003
$ws.Range("N80").Value = -4772.6365

$ws.Range("H83").Value = 628.4091
$ws.Range("I83").Value = 331.27274
$ws.Range("J83").Value = 925.5454999999999
$ws.Range("K83").Value = 2981.45466
$ws.Range("L83").Value = 8329.9095
$ws.Range("M83").Value = 2010.54534
$ws.Range("N83").Value = -18313.9095

$ws.Range("H100").Value = 3332.75
$ws.Range("I100").Value = 3537.5715
$ws.Range("J100").Value = 1899
$ws.Range("K100").Value = 3537.5715
$ws.Range("L100").Value = 1899
$ws.Range("M100").Value = -2996.5715
$ws.Range("N100").Value = -2981

$ws.Range("H137").Value = 4044.8333
$ws.Range("I137").Value = 3884.7144
$ws.Range("K137").Value = 11654.1432
$ws.Range("M137").Value = -9104.143199999999

$ws.Range("H141").Value = 1208.1666
$ws.Range("I141").Value = 1208.1666
$ws.Range("K141").Value = 3624.4998
$ws.Range("M141").Value = 1555.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10061
$ws.Range("I32").Value = 10063.48
$ws.Range("K32").Value = 10063.48
$ws.Range("M32").Value = -9776.48

$ws.Range("H63").Value = 19050.555
$ws.Range("I63").Value = 18779.285
$ws.Range("K63").Value = 18779.285
$ws.Range("M63").Value = -18093.285

$ws.Range("H66").Value = 19050.555
$ws.Range("I66").Value = 18779.285
$ws.Range("K66").Value = 93896.425
$ws.Range("M66").Value = -90464.425

$ws.Range("H74").Value = 2068.9333
$ws.Range("I74").Value = 2038.2858
$ws.Range("K74").Value = 2038.2858
$ws.Range("M74").Value = -1164.2858

$ws.Range("H77").Value = 2068.9333
$ws.Range("I77").Value = 2038.2858
$ws.Range("K77").Value = 10191.429
$ws.Range("M77").Value = -5823.429

$ws.Range("H131").Value = 66999
$ws.Range("J131").Value = 66999
$ws.Range("L131").Value = 66999
$ws.Range("N131").Value = -77079

$ws.Range("H133").Value = 59892.5
$ws.Range("J133").Value = 59892.5
$ws.Range("L133").Value = 59892.5
$ws.Range("N133").Value = -64952.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 369.6
$ws.Range("I22").Value = 377.55554
$ws.Range("J22").Value = 298
$ws.Range("K22").Value = 377.55554
$ws.Range("L22").Value = 298
$ws.Range("M22").Value = -204.55554
$ws.Range("N22").Value = -644

$ws.Range("H105").Value = 4458.231
$ws.Range("I105").Value = 4385.6665
$ws.Range("K105").Value = 4385.6665
$ws.Range("M105").Value = -2638.6665

$ws.Range("H107").Value = 890.5714
$ws.Range("I107").Value = 872.3333
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 872.3333
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1047.6667
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 7750
$ws.Range("J14").Value = 7750
$ws.Range("L14").Value = 7750
$ws.Range("N14").Value = -8090

$ws.Range("H22").Value = 3333761.2
$ws.Range("I22").Value = 502.25
$ws.Range("K22").Value = 502.25
$ws.Range("M22").Value = -152.25

$ws.Range("H31").Value = 2020.1666
$ws.Range("I31").Value = 2020.1666
$ws.Range("K31").Value = 2020.1666
$ws.Range("M31").Value = -1725.1666

$ws.Range("H34").Value = 2020.1666
$ws.Range("I34").Value = 2020.1666
$ws.Range("K34").Value = 2020.1666
$ws.Range("M34").Value = -1818.1666

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H105").Value = 2910.6667
$ws.Range("I105").Value = 2773.125
$ws.Range("K105").Value = 2773.125
$ws.Range("M105").Value = -1026.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 9010.125
$ws.Range("J94").Value = 9010.125
$ws.Range("L94").Value = 27030.375
$ws.Range("N94").Value = -28382.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2058.1177
$ws.Range("I122").Value = 1709.2142
$ws.Range("K122").Value = 5127.642599999999
$ws.Range("M122").Value = -2677.642599999999

$ws.Range("H128").Value = 85469.5
$ws.Range("J128").Value = 86365.39999999999
$ws.Range("L128").Value = 86365.39999999999
$ws.Range("N128").Value = -96325.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 6833
$ws.Range("I19").Value = 5000
$ws.Range("J19").Value = 7749.5
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 7749.5
$ws.Range("M19").Value = -4830
$ws.Range("N19").Value = -8089.5

$ws.Range("H46").Value = 2197.4
$ws.Range("I46").Value = 1774.4286
$ws.Range("J46").Value = 3184.3333
$ws.Range("K46").Value = 1774.4286
$ws.Range("L46").Value = 3184.3333
$ws.Range("M46").Value = -1586.4286
$ws.Range("N46").Value = -3560.3333

$ws.Range("H61").Value = 3038.6
$ws.Range("I61").Value = 3131.3333
$ws.Range("J61").Value = 2899.5
$ws.Range("K61").Value = 3131.3333
$ws.Range("L61").Value = 2899.5
$ws.Range("M61").Value = -2929.3333
$ws.Range("N61").Value = -3303.5

$ws.Range("H82").Value = 2184.111
$ws.Range("J82").Value = 2692.2
$ws.Range("L82").Value = 2692.2
$ws.Range("N82").Value = -3414.2

$ws.Range("H85").Value = 2184.111
$ws.Range("J85").Value = 2692.2
$ws.Range("L85").Value = 2692.2
$ws.Range("N85").Value = -5188.2

$ws.Range("H113").Value = 3038.6
$ws.Range("I113").Value = 3131.3333
$ws.Range("J113").Value = 2899.5
$ws.Range("K113").Value = 3131.3333
$ws.Range("L113").Value = 2899.5
$ws.Range("M113").Value = -961.3332999999998
$ws.Range("N113").Value = -7239.5

$ws.Range("H130").Value = 66660
$ws.Range("J130").Value = 66660
$ws.Range("L130").Value = 66660
$ws.Range("N130").Value = -76700

$ws.Range("H136").Value = 3090
$ws.Range("I136").Value = 2544.4285
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 7633.2855
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -5083.2855
$ws.Range("N136").Value = -20098.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4644.636
$ws.Range("I81").Value = 2636.375
$ws.Range("K81").Value = 5272.75
$ws.Range("M81").Value = -4211.75

$ws.Range("H84").Value = 4644.636
$ws.Range("I84").Value = 2636.375
$ws.Range("K84").Value = 26363.75
$ws.Range("M84").Value = -21059.75

$ws.Range("H130").Value = 26999
$ws.Range("J130").Value = 26999
$ws.Range("L130").Value = 26999
$ws.Range("N130").Value = -37039

